# Registration "required fields" regression-test data: a new worksheet is
# added after "ContactUs" containing eight rows of otherwise-identical
# registration records, each missing exactly one required field so the
# automation suite can assert the right validation message fires.

$wb = $excel.ActiveWorkbook

# Add the new sheet at the very end of the tab strip (after ContactUs).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "RegistrationRequiredFields"

# Column layout: name, email, password, day, month, year, firstName,
# lastName, address1, country, state, city, zipcode, mobile
$name     = "TempUser"
$email    = "reg_req_{{unique}}@ex.com"
$password = "Test1234"
$day      = "1"
$month    = "January"
$year     = "2000"
$first    = "John"
$last     = "Doe"
$address1 = "123 Street"
$country  = "Canada"
$state    = "ON"
$city     = "Toronto"
$zipcode  = "M1A1A1"
$mobile   = "1234567"

$fields = @($name, $email, $password, $day, $month, $year, $first, $last, $address1, $country, $state, $city, $zipcode, $mobile)

for ($row = 1; $row -le 8; $row++) {
    $missingCol = $row + 2   # row1 blanks col3 (password), row2 blanks col7 (firstName), ...
    if ($row -eq 1) { $missingCol = 3 }
    elseif ($row -eq 2) { $missingCol = 7 }
    elseif ($row -eq 3) { $missingCol = 9 }
    elseif ($row -eq 4) { $missingCol = 10 }
    elseif ($row -eq 5) { $missingCol = 11 }
    elseif ($row -eq 6) { $missingCol = 12 }
    elseif ($row -eq 7) { $missingCol = 13 }
    elseif ($row -eq 8) { $missingCol = 14 }

    for ($col = 1; $col -le 14; $col++) {
        if ($col -eq $missingCol) {
            $ws.Cells.Item($row, $col).Value = ""
        } else {
            $ws.Cells.Item($row, $col).Value = $fields[$col - 1]
        }
    }
}

$ws.Range("K5").Select()
